$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 495.11765
$ws.Range("I39").Value = 114.333336
$ws.Range("J39").Value = 923.5
$ws.Range("K39").Value = 343.000008
$ws.Range("L39").Value = 2770.5
$ws.Range("M39").Value = -47.00000799999998
$ws.Range("N39").Value = -3362.5
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()
$ws.Range("H135").Value = 5102786
$ws.Range("I135").Value = 5952922
$ws.Range("J135").Value = 1970.2858
$ws.Range("K135").Value = 53576298
$ws.Range("L135").Value = 17732.5722
$ws.Range("M135").Value = -53573763
$ws.Range("N135").Value = -22802.5722

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1439.2222
$ws.Range("I2").Value = 1231.3125
$ws.Range("J2").Value = 1741.6364
$ws.Range("K2").Value = 1231.3125
$ws.Range("L2").Value = 1741.6364
$ws.Range("M2").Value = -1118.3125
$ws.Range("N2").Value = -1967.6364
$ws.Range("H33").Value = 15000
$ws.Range("J33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15658
$ws.Range("H47").Value = 16593.334
$ws.Range("J47").Value = 16593.334
$ws.Range("L47").Value = 16593.334
$ws.Range("N47").Value = -18043.334
$ws.Range("H49").Value = 13000
$ws.Range("J49").Value = 13000
$ws.Range("L49").Value = 13000
$ws.Range("N49").Value = -13520
$ws.Range("H104").Value = 30225
$ws.Range("J104").Value = 30225
$ws.Range("L104").Value = 30225
$ws.Range("N104").Value = -37213
$ws.Range("H116").Value = 1439.2222
$ws.Range("I116").Value = 1231.3125
$ws.Range("J116").Value = 1741.6364
$ws.Range("K116").Value = 1231.3125
$ws.Range("L116").Value = 1741.6364
$ws.Range("M116").Value = 1062.6875
$ws.Range("N116").Value = -6329.6364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1439.2222
$ws.Range("I3").Value = 1231.3125
$ws.Range("J3").Value = 1741.6364
$ws.Range("K3").Value = 1231.3125
$ws.Range("L3").Value = 1741.6364
$ws.Range("M3").Value = -1117.3125
$ws.Range("N3").Value = -1969.6364
$ws.Range("H105").Value = 71430030
$ws.Range("I105").Value = 1383.5
$ws.Range("J105").Value = 166668240
$ws.Range("K105").Value = 1383.5
$ws.Range("L105").Value = 166668240
$ws.Range("M105").Value = 363.5
$ws.Range("N105").Value = -166671734

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 15718
$ws.Range("I33").Value = 2154
$ws.Range("J33").Value = 22500
$ws.Range("K33").Value = 2154
$ws.Range("L33").Value = 22500
$ws.Range("M33").Value = -1775
$ws.Range("N33").Value = -23258
$ws.Range("H57").Value = 6000
$ws.Range("J57").Value = 6000
$ws.Range("L57").Value = 6000
$ws.Range("N57").Value = -7120
$ws.Range("H74").Value = 20555.334
$ws.Range("J74").Value = 18333
$ws.Range("L74").Value = 18333
$ws.Range("N74").Value = -20081
$ws.Range("H77").Value = 20555.334
$ws.Range("J77").Value = 18333
$ws.Range("L77").Value = 54999
$ws.Range("N77").Value = -63735
$ws.Range("H86").Value = 27810224
$ws.Range("I86").Value = 71433640
$ws.Range("J86").Value = 49864.363
$ws.Range("K86").Value = 71433640
$ws.Range("L86").Value = 49864.363
$ws.Range("M86").Value = -71432517
$ws.Range("N86").Value = -52110.363
$ws.Range("H89").Value = 27810224
$ws.Range("I89").Value = 71433640
$ws.Range("J89").Value = 49864.363
$ws.Range("K89").Value = 357168200
$ws.Range("L89").Value = 249321.815
$ws.Range("M89").Value = -357162584
$ws.Range("N89").Value = -260553.815

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 418.7143
$ws.Range("I41").Value = 170
$ws.Range("J41").Value = 750.3333
$ws.Range("K41").Value = 510
$ws.Range("L41").Value = 2250.9999
$ws.Range("M41").Value = -172
$ws.Range("N41").Value = -2926.9999
$ws.Range("H110").Value = 2027
$ws.Range("I110").Value = 2027
$ws.Range("K110").Value = 6081
$ws.Range("M110").Value = -1991

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 6000
$ws.Range("J49").Value = 6000
$ws.Range("L49").Value = 6000
$ws.Range("N49").Value = -6368

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 10000
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10820
$ws.Range("H47").Value = 10000
$ws.Range("J47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -10980
$ws.Range("H48").Value = 13260
$ws.Range("I48").Value = 10000
$ws.Range("K48").Value = 10000
$ws.Range("M48").Value = -9339
$ws.Range("H52").Value = 10000
$ws.Range("J52").Value = 10000
$ws.Range("L52").Value = 10000
$ws.Range("N52").Value = -10466
$ws.Range("H53").Value = 13738
$ws.Range("J53").Value = 13738
$ws.Range("L53").Value = 13738
$ws.Range("N53").Value = -14774
$ws.Range("H62").Value = 40065
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 40065
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 40065
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -41313
$ws.Range("H65").Value = 40065
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 40065
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 120195
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -126435
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 15450
$ws.Range("J47").Value = 15450
$ws.Range("L47").Value = 15450
$ws.Range("N47").Value = -16594
$ws.Range("H48").Value = 5000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 5000
$ws.Range("N48").Value = -6138
